# OLX Monitor 2026-02-15 21:59
# Adds a new monitoring snapshot (two check rows: 21:58 and 21:59) to each
# profile sheet, and bumps the "last checked" timestamp on the summary
# sheet to 21:59.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) PODSUMOWANIE (summary) sheet: "Data ostatniego sprawdzenia" column B
#    moves from 21:51 -> 21:59 for every profile row.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("PODSUMOWANIE")
$summary.Range("B2").Value = "2026-02-15 21:59"
$summary.Range("B3").Value = "2026-02-15 21:59"
$summary.Range("B4").Value = "2026-02-15 21:59"
$summary.Range("B5").Value = "2026-02-15 21:59"
$summary.Range("B6").Value = "2026-02-15 21:59"

# ---------------------------------------------------------------------
# Helper routine (inlined per-sheet below, PS COM has no closures over
# $ws reliably across calls in this host, so we repeat the steps).
# Each profile sheet gets two new rows appended at the bottom:
#   row 8 (21:58 check, "even" shaded style like row 6)
#   row 9 (21:59 check, "odd" plain style like row 7)
# Row formats are cloned from the last two existing rows (6 and 7) via
# PasteSpecial(formats) so the alternating-stripe styling is preserved
# exactly, then values are written on top.
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# 2) wszystkie-lublin  (totals stay at 432, nothing new/removed)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("wszystkie-lublin")

$ws.Range("A6:I6").Copy()
$ws.Range("A8:I8").PasteSpecial(-4122)
$ws.Rows.Item(8).RowHeight = 18
$ws.Range("I8").ClearContents()
$ws.Range("A8").Value = "2026-02-15 21:58"
$ws.Range("B8").Value = 432
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = "—"
$ws.Range("G8").Value = "—"
$ws.Range("H8").Value = "OK"

$ws.Range("A7:I7").Copy()
$ws.Range("A9:I9").PasteSpecial(-4122)
$ws.Rows.Item(9).RowHeight = 18
$ws.Range("I9").ClearContents()
$ws.Range("A9").Value = "2026-02-15 21:59"
$ws.Range("B9").Value = 432
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = "—"
$ws.Range("G9").Value = "—"
$ws.Range("H9").Value = "OK"

# ---------------------------------------------------------------------
# 3) artymiuk (totals stay at 0)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("artymiuk")

$ws.Range("A6:I6").Copy()
$ws.Range("A8:I8").PasteSpecial(-4122)
$ws.Rows.Item(8).RowHeight = 18
$ws.Range("I8").ClearContents()
$ws.Range("A8").Value = "2026-02-15 21:58"
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = "—"
$ws.Range("G8").Value = "—"
$ws.Range("H8").Value = "OK"

$ws.Range("A7:I7").Copy()
$ws.Range("A9:I9").PasteSpecial(-4122)
$ws.Rows.Item(9).RowHeight = 18
$ws.Range("I9").ClearContents()
$ws.Range("A9").Value = "2026-02-15 21:59"
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = "—"
$ws.Range("G9").Value = "—"
$ws.Range("H9").Value = "OK"

# ---------------------------------------------------------------------
# 4) poqui (totals stay at 5, "Szczegóły" id list carried to column I)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("poqui")

$ws.Range("A6:I6").Copy()
$ws.Range("A8:I8").PasteSpecial(-4122)
$ws.Rows.Item(8).RowHeight = 18
$ws.Range("A8").Value = "2026-02-15 21:58"
$ws.Range("B8").Value = 5
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = "—"
$ws.Range("G8").Value = "—"
$ws.Range("H8").Value = "OK"
$ws.Range("I8").Value = "18KAEc|17vbYq|1951OR|17NeTz|183ger"

$ws.Range("A7:I7").Copy()
$ws.Range("A9:I9").PasteSpecial(-4122)
$ws.Rows.Item(9).RowHeight = 18
$ws.Range("A9").Value = "2026-02-15 21:59"
$ws.Range("B9").Value = 5
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = "—"
$ws.Range("G9").Value = "—"
$ws.Range("H9").Value = "OK"
$ws.Range("I9").Value = "1951OR|17NeTz|17vbYq|183ger|18KAEc"

# ---------------------------------------------------------------------
# 5) stylowepokoje (totals stay at 2, "Szczegóły" id list carried to I)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("stylowepokoje")

$ws.Range("A6:I6").Copy()
$ws.Range("A8:I8").PasteSpecial(-4122)
$ws.Rows.Item(8).RowHeight = 18
$ws.Range("A8").Value = "2026-02-15 21:58"
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = "—"
$ws.Range("G8").Value = "—"
$ws.Range("H8").Value = "OK"
$ws.Range("I8").Value = "16ZeYm|195dLc"

$ws.Range("A7:I7").Copy()
$ws.Range("A9:I9").PasteSpecial(-4122)
$ws.Rows.Item(9).RowHeight = 18
$ws.Range("A9").Value = "2026-02-15 21:59"
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = "—"
$ws.Range("G9").Value = "—"
$ws.Range("H9").Value = "OK"
$ws.Range("I9").Value = "16ZeYm|195dLc"

# ---------------------------------------------------------------------
# 6) villahome (totals stay at 0)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("villahome")

$ws.Range("A6:I6").Copy()
$ws.Range("A8:I8").PasteSpecial(-4122)
$ws.Rows.Item(8).RowHeight = 18
$ws.Range("I8").ClearContents()
$ws.Range("A8").Value = "2026-02-15 21:58"
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = "—"
$ws.Range("G8").Value = "—"
$ws.Range("H8").Value = "OK"

$ws.Range("A7:I7").Copy()
$ws.Range("A9:I9").PasteSpecial(-4122)
$ws.Rows.Item(9).RowHeight = 18
$ws.Range("I9").ClearContents()
$ws.Range("A9").Value = "2026-02-15 21:59"
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = "—"
$ws.Range("G9").Value = "—"
$ws.Range("H9").Value = "OK"
